$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repair")

# Widen column A to fit the longer labels being added (closest value this
# engine's character-grid rounding can reach to the target 93.6).
$ws.Columns.Item(1).ColumnWidth = 92.8

# Row 29 used to be "Mean Repair Time Relative to Surface Nets Mesh Generation
# Time" (the last row before "Success Rate"). Insert 6 blank rows right after
# it so the sheet grows from 30 rows to 36 rows; this pushes the existing
# "Success Rate" row from 30 down to 36.
$ws.Range("A30:A35").EntireRow.Insert()

# Copy the label/value formatting (bold+border label, percentage value) down
# into the newly inserted rows so they match the rest of the block.
$ws.Range("A28:B28").Copy()
$ws.Range("A30:B35").PasteSpecial(-4122)
$ws.Range("B29").NumberFormat = $ws.Range("B28").NumberFormat
$ws.Range("B36").NumberFormat = $ws.Range("B28").NumberFormat

# Re-label row 29 and fill in the new Cuberille stats rows (30-31).
$ws.Range("A29").Value = "Min Repair Time Relative to Cuberille Mesh Generation Time"
$ws.Range("B29").Value = 0.006947098007963789

$ws.Range("A30").Value = "Max Repair Time Relative to Cuberille Mesh Generation Time"
$ws.Range("B30").Value = 2.235620867295733

$ws.Range("A31").Value = "Standard Deviation Repair Time Relative to Cuberille Mesh Generation Time"
$ws.Range("B31").Value = 0.2488410955062405

# The original "Mean ... Surface Nets ..." row now lives at 32, followed by
# its new Min/Max/StdDev siblings (33-35).
$ws.Range("A32").Value = "Mean Repair Time Relative to Surface Nets Mesh Generation Time"
$ws.Range("B32").Value = 0.04746570377831323

$ws.Range("A33").Value = "Min Repair Time Relative to Surface Nets Mesh Generation Time"
$ws.Range("B33").Value = 0.002070496253869231

$ws.Range("A34").Value = "Max Repair Time Relative to Surface Nets Mesh Generation Time"
$ws.Range("B34").Value = 0.9495387351385597

$ws.Range("A35").Value = "Standard Deviation Repair Time Relative to Surface Nets Mesh Generation Time"
$ws.Range("B35").Value = 0.07246700117742484

# Row 36 already holds the shifted-down "Success Rate" / 1 from the insert.
